$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F width
$ws.Columns.Item(6).ColumnWidth = 19.375

# Row 2 - header, green fill (same style as E2)
$ws.Range("F2").Value = "Tuần 3"
$ws.Range("F2").Interior.Color = $ws.Range("E2").Interior.Color

# Row 3 - plain values
$ws.Range("F3").Value = "Validate các màn hình 4/"

# Row 5 - plain values (entered before row 4, matching shared-string order)
$ws.Range("F5").Value = "Validate các màn hình 2/ + 5/"

# Row 4 - plain values
$ws.Range("F4").Value = "Validate các màn hình 1/ + 3/"

# Row 6 - F6 inherits the style E6 currently has (red font + yellow fill)
$ws.Range("F6").Value = "Thứ 4, 10/4/2019"
$ws.Range("F6").Font.Color = $ws.Range("E6").Font.Color
$ws.Range("F6").Interior.Color = $ws.Range("E6").Interior.Color

# E6 now gets a new style: red font + white (theme background 1) fill
$ws.Range("E6").Interior.ThemeColor = 2

# Update the view: top-left cell and selection
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F7").Select()
